$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D, E and G hold numeric-looking values ("305.69", "0.89%", "12") that
# are stored as literal TEXT in the source workbook (t="inlineStr"). Force the
# cells to Text format first so Excel does not silently convert the new values
# assigned below into real numbers/percentages.
$deRng = $ws.Range("D2:E51")
$gRng = $ws.Range("G2:G51")
$deRng.NumberFormat = "@"
$gRng.NumberFormat = "@"

$ws.Range('D2').Value = '306.01'
$ws.Range('E2').Value = '0.86%'
$ws.Range('G2').Value = '13'
$ws.Range('D3').Value = '35.98'
$ws.Range('E3').Value = '-3.91%'
$ws.Range('G3').Value = '13'
$ws.Range('D4').Value = '5.089'
$ws.Range('E4').Value = '1.52%'
$ws.Range('G4').Value = '13'
$ws.Range('D5').Value = '0.07977'
$ws.Range('E5').Value = '1.41%'
$ws.Range('G5').Value = '13'
$ws.Range('E6').Value = '-4.44%'
$ws.Range('G6').Value = '13'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '4.154'
$ws.Range('E7').Value = '3.16%'
$ws.Range('G7').Value = '13'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D8').Value = '7.971'
$ws.Range('E8').Value = '-0.52%'
$ws.Range('G8').Value = '13'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '0.9228'
$ws.Range('E9').Value = '0.61%'
$ws.Range('G9').Value = '13'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '0.09661'
$ws.Range('E10').Value = '1.19%'
$ws.Range('G10').Value = '13'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1852'
$ws.Range('E11').Value = '-1.13%'
$ws.Range('G11').Value = '13'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.08736'
$ws.Range('E12').Value = '1.43%'
$ws.Range('G12').Value = '13'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.03561'
$ws.Range('E13').Value = '-1.02%'
$ws.Range('G13').Value = '13'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.09955'
$ws.Range('E14').Value = '0.01%'
$ws.Range('G14').Value = '13'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001431'
$ws.Range('E15').Value = '-3.52%'
$ws.Range('G15').Value = '13'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005710'
$ws.Range('E16').Value = '0.29%'
$ws.Range('G16').Value = '13'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.473'
$ws.Range('E17').Value = '0.35%'
$ws.Range('G17').Value = '13'
$ws.Range('D18').Value = '2.666'
$ws.Range('E18').Value = '18.54%'
$ws.Range('G18').Value = '13'
$ws.Range('D19').Value = '0.3387'
$ws.Range('E19').Value = '0.20%'
$ws.Range('G19').Value = '13'
$ws.Range('E20').Value = '2.42%'
$ws.Range('G20').Value = '13'
$ws.Range('D21').Value = '5.161'
$ws.Range('E21').Value = '8.82%'
$ws.Range('G21').Value = '13'
$ws.Range('D22').Value = '0.2207'
$ws.Range('E22').Value = '0.35%'
$ws.Range('G22').Value = '13'
$ws.Range('D23').Value = '0.04560'
$ws.Range('E23').Value = '-0.40%'
$ws.Range('G23').Value = '13'
$ws.Range('E24').Value = '0.59%'
$ws.Range('G24').Value = '13'
$ws.Range('D25').Value = '0.004907'
$ws.Range('E25').Value = '2.41%'
$ws.Range('G25').Value = '13'
$ws.Range('D26').Value = '0.0001300'
$ws.Range('E26').Value = '-6.99%'
$ws.Range('G26').Value = '13'
$ws.Range('D27').Value = '0.0004749'
$ws.Range('E27').Value = '-0.03%'
$ws.Range('G27').Value = '13'
$ws.Range('G28').Value = '13'
$ws.Range('G29').Value = '13'
$ws.Range('G30').Value = '13'
$ws.Range('G31').Value = '13'
$ws.Range('G32').Value = '13'
$ws.Range('G33').Value = '13'
$ws.Range('G34').Value = '13'
$ws.Range('G35').Value = '13'
$ws.Range('G36').Value = '13'
$ws.Range('G37').Value = '13'
$ws.Range('G38').Value = '13'
$ws.Range('D39').Value = '0.01861'
$ws.Range('E39').Value = '2.60%'
$ws.Range('G39').Value = '13'
$ws.Range('D40').Value = '0.04759'
$ws.Range('E40').Value = '0.61%'
$ws.Range('G40').Value = '13'
$ws.Range('D41').Value = '0.007879'
$ws.Range('E41').Value = '-3.19%'
$ws.Range('G41').Value = '13'
$ws.Range('D42').Value = '0.1400'
$ws.Range('E42').Value = '0.34%'
$ws.Range('G42').Value = '13'
$ws.Range('D43').Value = '0.007756'
$ws.Range('E43').Value = '2.62%'
$ws.Range('G43').Value = '13'
$ws.Range('D44').Value = '0.002205'
$ws.Range('E44').Value = '-0.52%'
$ws.Range('G44').Value = '13'
$ws.Range('D45').Value = '0.01127'
$ws.Range('E45').Value = '7.73%'
$ws.Range('G45').Value = '13'
$ws.Range('D46').Value = '0.00006259'
$ws.Range('E46').Value = '1.60%'
$ws.Range('G46').Value = '13'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').Value = '0.17%'
$ws.Range('G47').Value = '13'
$ws.Range('E48').Value = '0.34%'
$ws.Range('G48').Value = '13'
$ws.Range('D49').Value = '49.59'
$ws.Range('E49').Value = '24.68%'
$ws.Range('G49').Value = '13'
$ws.Range('D50').Value = '0.002000'
$ws.Range('E50').Value = '-25.66%'
$ws.Range('G50').Value = '13'
$ws.Range('D51').Value = '0.00002100'
$ws.Range('E51').Value = '0.17%'
$ws.Range('G51').Value = '13'

# Reset cell style back to Normal so no stray "Text" number-format style is
# left referenced on these cells (matches the original, unstyled cells).
$deRng.Style = "Normal"
$gRng.Style = "Normal"
